$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data range A2:M64 ascending by column A (header row excluded),
# mirroring a manual Data > Sort operation performed in the Excel UI.
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A2:A64"))
$sort.SetRange($ws.Range("A2:M64"))
$sort.Header = 2
$sort.Apply()

# Update the active selection to match the post-sort cursor position.
$ws.Range("E20").Select()
